# Update F-column ('想去人数' / interest counts) values across all four sheets
# of the workbook, per the source diff.
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 35
$ws.Range("F4").Value = 842
$ws.Range("F6").Value = 4
$ws.Range("F7").Value = 9618
$ws.Range("F8").Value = 39
$ws.Range("F10").Value = 2159
$ws.Range("F11").Value = 50
$ws.Range("F12").Value = 1643
$ws.Range("F13").Value = 2766
$ws.Range("F14").Value = 142
$ws.Range("F15").Value = 4119
$ws.Range("F16").Value = 345
$ws.Range("F17").Value = 169
$ws.Range("F18").Value = 133
$ws.Range("F19").Value = 523
$ws.Range("F21").Value = 33
$ws.Range("F23").Value = 85
$ws.Range("F25").Value = 3963
$ws.Range("F27").Value = 3446
$ws.Range("F28").Value = 1106
$ws.Range("F29").Value = 203
$ws.Range("F30").Value = 505
$ws.Range("F33").Value = 319
$ws.Range("F34").Value = 422
$ws.Range("F35").Value = 311

# Sheet 2: 演出 (Performances)
$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 2
$ws.Range("F3").Value = 24

# Sheet 3: 本地生活 (Local Life)
$ws = $wb.Worksheets.Item(3)
$ws.Range("F3").Value = 1009

# Sheet 4: 全部类型 (All Types)
$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 35
$ws.Range("F4").Value = 1009
$ws.Range("F6").Value = 842
$ws.Range("F8").Value = 4
$ws.Range("F9").Value = 9618
$ws.Range("F10").Value = 39
$ws.Range("F12").Value = 2160
$ws.Range("F13").Value = 50
$ws.Range("F14").Value = 1643
$ws.Range("F15").Value = 2
$ws.Range("F16").Value = 2766
$ws.Range("F17").Value = 142
$ws.Range("F18").Value = 4120
$ws.Range("F19").Value = 345
$ws.Range("F20").Value = 169
$ws.Range("F21").Value = 133
$ws.Range("F22").Value = 523
$ws.Range("F24").Value = 33
$ws.Range("F25").Value = 24
$ws.Range("F27").Value = 85
$ws.Range("F29").Value = 3963
$ws.Range("F31").Value = 3446
$ws.Range("F32").Value = 1106
$ws.Range("F33").Value = 203
$ws.Range("F34").Value = 505
$ws.Range("F37").Value = 319
$ws.Range("F38").Value = 422
$ws.Range("F39").Value = 311

Write-Output "Done updating interest counts."
